$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.942.71"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "3.768.37"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "426.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.617"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.78%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.724"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("E10").Value = "  -11.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000301"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -16.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.36%  "
$ws.Range("D14").Value = "4.370.30"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.137"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "3.782.41"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("E19").Value = "  +4.97%  "
$ws.Range("D20").Value = "66.060.82"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "402.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.65%  "
$ws.Range("E23").Value = "  +6.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "36.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +33.48%  "
$ws.Range("E27").Value = "  +5.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.51%  "
$ws.Range("E29").Value = "  -6.14%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.99%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "699.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.25%  "
$ws.Range("E32").Value = "  +11.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.69"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.73%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +33.85%  "
$ws.Range("E37").Value = "  -3.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "56.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E39").Value = "  +2.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +32.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.140"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.94%  "
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("D44").Value = "0.0₃0653"
$ws.Range("E44").Value = "  -11.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.319"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.17%  "
$ws.Range("E48").Value = "  +3.87%  "
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "138.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.83%  "
